$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 109
$ws.Cells.Item(4, 6).Value = 1823
$ws.Cells.Item(5, 6).Value = 143
$ws.Cells.Item(6, 6).Value = 3844
$ws.Cells.Item(7, 6).Value = 511
$ws.Cells.Item(8, 6).Value = 1033
$ws.Cells.Item(9, 6).Value = 1286
$ws.Cells.Item(10, 6).Value = 646
$ws.Cells.Item(11, 6).Value = 360
$ws.Cells.Item(13, 6).Value = 2133
$ws.Cells.Item(14, 6).Value = 385
$ws.Cells.Item(15, 6).Value = 637849
$ws.Cells.Item(16, 6).Value = 1567
$ws.Cells.Item(17, 6).Value = 448
$ws.Cells.Item(18, 6).Value = 1380
$ws.Cells.Item(19, 6).Value = 658
$ws.Cells.Item(21, 6).Value = 1228
$ws.Cells.Item(22, 6).Value = 2105
$ws.Cells.Item(23, 6).Value = 1083
$ws.Cells.Item(24, 6).Value = 2636
$ws.Cells.Item(25, 6).Value = 1506
$ws.Cells.Item(26, 6).Value = 722
$ws.Cells.Item(27, 6).Value = 1472
$ws.Cells.Item(28, 6).Value = 20
$ws.Cells.Item(29, 6).Value = 513
$ws.Cells.Item(30, 6).Value = 1060
$ws.Cells.Item(31, 6).Value = 221
$ws.Cells.Item(32, 6).Value = 1060
$ws.Cells.Item(34, 6).Value = 66
$ws.Cells.Item(35, 6).Value = 1974
$ws.Cells.Item(36, 6).Value = 1286
$ws.Cells.Item(37, 6).Value = 549
$ws.Cells.Item(38, 6).Value = 607
$ws.Cells.Item(39, 6).Value = 1118
$ws.Cells.Item(41, 6).Value = 189
$ws.Cells.Item(43, 6).Value = 2509
$ws.Cells.Item(44, 6).Value = 196
$ws.Cells.Item(45, 6).Value = 956
$ws.Cells.Item(46, 6).Value = 3043
$ws.Cells.Item(49, 6).Value = 663
$ws.Cells.Item(50, 6).Value = 37

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 61
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 461
$ws.Cells.Item(10, 6).Value = 88
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 144079
$ws.Cells.Item(17, 6).Value = 91
$ws.Cells.Item(19, 6).Value = 324
$ws.Cells.Item(21, 6).Value = 391
$ws.Cells.Item(22, 6).Value = 391
$ws.Cells.Item(23, 6).Value = 96
$ws.Cells.Item(25, 6).Value = 96
$ws.Cells.Item(26, 6).Value = 84
$ws.Cells.Item(27, 6).Value = 492
$ws.Cells.Item(28, 6).Value = 87
$ws.Cells.Item(32, 6).Value = 286
$ws.Cells.Item(33, 6).Value = 257
$ws.Cells.Item(35, 6).Value = 48
$ws.Cells.Item(36, 6).Value = 19
$ws.Cells.Item(38, 6).Value = 107
$ws.Cells.Item(40, 6).Value = 180
$ws.Cells.Item(45, 6).Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3098
$ws.Cells.Item(5, 6).Value = 4855
$ws.Cells.Item(8, 6).Value = 800
$ws.Cells.Item(9, 6).Value = 1113
$ws.Cells.Item(10, 6).Value = 614
$ws.Cells.Item(11, 6).Value = 1548
$ws.Cells.Item(13, 6).Value = 1747

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 800
$ws.Cells.Item(3, 6).Value = 614
$ws.Cells.Item(5, 6).Value = 109
$ws.Cells.Item(6, 6).Value = 1823
$ws.Cells.Item(7, 6).Value = 1747
$ws.Cells.Item(8, 6).Value = 3844
$ws.Cells.Item(10, 6).Value = 511
$ws.Cells.Item(11, 6).Value = 1286
$ws.Cells.Item(12, 6).Value = 646
$ws.Cells.Item(13, 6).Value = 360
$ws.Cells.Item(14, 6).Value = 2133
$ws.Cells.Item(15, 6).Value = 385
$ws.Cells.Item(16, 6).Value = 637860
$ws.Cells.Item(17, 6).Value = 461
$ws.Cells.Item(18, 6).Value = 88
$ws.Cells.Item(19, 6).Value = 1567
$ws.Cells.Item(20, 6).Value = 144079
$ws.Cells.Item(21, 6).Value = 448
$ws.Cells.Item(22, 6).Value = 1380
$ws.Cells.Item(23, 6).Value = 658
$ws.Cells.Item(25, 6).Value = 1228
$ws.Cells.Item(26, 6).Value = 2105
$ws.Cells.Item(27, 6).Value = 1083
$ws.Cells.Item(28, 6).Value = 2636
$ws.Cells.Item(29, 6).Value = 1506
$ws.Cells.Item(30, 6).Value = 722
$ws.Cells.Item(32, 6).Value = 1472
$ws.Cells.Item(33, 6).Value = 391
$ws.Cells.Item(34, 6).Value = 513
$ws.Cells.Item(35, 6).Value = 96
$ws.Cells.Item(36, 6).Value = 1060
$ws.Cells.Item(37, 6).Value = 1060
$ws.Cells.Item(39, 6).Value = 66
$ws.Cells.Item(40, 6).Value = 1974
$ws.Cells.Item(41, 6).Value = 1286
$ws.Cells.Item(42, 6).Value = 549
$ws.Cells.Item(43, 6).Value = 612
$ws.Cells.Item(44, 6).Value = 1118
$ws.Cells.Item(45, 6).Value = 286
$ws.Cells.Item(46, 6).Value = 286
$ws.Cells.Item(47, 6).Value = 257
$ws.Cells.Item(48, 6).Value = 2509
$ws.Cells.Item(49, 6).Value = 196
$ws.Cells.Item(50, 6).Value = 956
$ws.Cells.Item(51, 6).Value = 3043
$ws.Cells.Item(52, 6).Value = 663
